# Diary workbook update: add a new diary entry (row 30) for "14 marras".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new entry in row 30 -------------------------------------------------
$ws.Range("A30").Value = "14 marras"
$ws.Range("B30").Value = "8.45-11.15, 12.00-13.00"
$ws.Range("C30").Value = "Kovien kappaleiden demon debug, tsemppi"
$ws.Range("G30").Value = 3.5

# Match the formatting used by the other rows in the table:
#  - column B uses the "time range" style (no wrap)
#  - column C uses the wrapped-text style
$ws.Range("B9").Copy()
$ws.Range("B30").PasteSpecial(-4122)

$ws.Range("C29").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row height to match the other wrapped rows in the sheet
$ws.Rows.Item(30).RowHeight = 29

# --- Update selection state --------------------------------------------------
$null = $ws.Range("H4").Select()

# --- Recalculate so the SUM() in H3 picks up the new hours -------------------
$excel.Calculate()
